$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TransportationPlan")

# Insert a new column before the current column B (ProductName), shifting
# ProductName/PlantName/CustomerName/X_transportation/X_transportation_Solution
# one column to the right.
$ws.Columns.Item(2).Insert()

# The inserted column picks up column A's formatting (border/bold style) for
# every row; strip that from the data rows so only the header keeps it.
$ws.Range("B2:B9").ClearFormats()

# Header for the new column - reuse the existing header style (bold/bordered)
# from the neighboring header cell instead of building a brand new style.
$ws.Range("B1").Value = "index"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the new "index" column with the same numbering as column A
for ($i = 0; $i -le 7; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $i
}
